# Updates the betting-odds figures on "Sheet1" to match the
# 2024-11-04 FlashScore refresh: 97 numeric cells across rows
# 2, 3, 4, 5, 8, 9 and 13 change value; nothing else in the
# workbook (text columns, headers, other rows) is touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Fulham vs Brentford)
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 4.33
$ws.Range("Q2").Value = 1.73
$ws.Range("R2").Value = 2.1

# Row 3 (Lazio vs Cagliari)
$ws.Range("AW3").Value = 151

# Row 4 (Celta Vigo vs Getafe)
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 4.1
$ws.Range("L4").Value = 5
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.3
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 1.4
$ws.Range("S4").Value = 1.62
$ws.Range("T4").Value = 2.2
$ws.Range("AH4").Value = 8
$ws.Range("AT4").Value = 2.2
$ws.Range("BC4").Value = 451

# Row 5 (Atl. Tucuman vs Sarmiento Junin)
$ws.Range("G5").Value = 1.75
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.5
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 5.5
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 1.5
$ws.Range("T5").Value = 2.5
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 7
$ws.Range("Z5").Value = 13
$ws.Range("AB5").Value = 34
$ws.Range("AD5").Value = 6.5
$ws.Range("AE5").Value = 21
$ws.Range("AH5").Value = 10
$ws.Range("AI5").Value = 23
$ws.Range("AJ5").Value = 17
$ws.Range("AK5").Value = 51
$ws.Range("AL5").Value = 41
$ws.Range("AM5").Value = 51
$ws.Range("AN5").Value = 3.5
$ws.Range("AO5").Value = 9.5
$ws.Range("AQ5").Value = 34
$ws.Range("AT5").Value = 2.5
$ws.Range("AU5").Value = 9.5
$ws.Range("AX5").Value = 6.5
$ws.Range("AY5").Value = 29
$ws.Range("BA5").Value = 126
$ws.Range("BB5").Value = 151
$ws.Range("BC5").Value = 401

# Row 8 (Tecnico U. vs U. Catolica)
$ws.Range("G8").Value = 3.7
$ws.Range("I8").Value = 1.95
$ws.Range("J8").Value = 4.5
$ws.Range("L8").Value = 2.63
$ws.Range("N8").Value = 9
$ws.Range("W8").Value = 9.5
$ws.Range("X8").Value = 19
$ws.Range("AI8").Value = 8.5
$ws.Range("AL8").Value = 17
$ws.Range("AO8").Value = 23
$ws.Range("AQ8").Value = 81
$ws.Range("AV8").Value = 67

# Row 9 (AC Ajaccio vs Metz)
$ws.Range("G9").Value = 4
$ws.Range("I9").Value = 2
$ws.Range("S9").Value = 1.53
$ws.Range("T9").Value = 2.38
$ws.Range("W9").Value = 9
$ws.Range("X9").Value = 19
$ws.Range("AB9").Value = 41
$ws.Range("AI9").Value = 8.5
$ws.Range("AJ9").Value = 9.5
$ws.Range("AN9").Value = 5.5
$ws.Range("AS9").Value = 301
$ws.Range("AT9").Value = 2.38

# Row 13 (Nacional vs Santa Clara)
$ws.Range("H13").Value = 3.1
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 3.4
$ws.Range("M13").Value = 1.08
$ws.Range("N13").Value = 8
$ws.Range("O13").Value = 1.4
$ws.Range("P13").Value = 2.75
$ws.Range("Q13").Value = 2.35
$ws.Range("R13").Value = 1.57
$ws.Range("S13").Value = 1.5
$ws.Range("T13").Value = 2.5
$ws.Range("U13").Value = 1.95
$ws.Range("V13").Value = 1.8
$ws.Range("W13").Value = 7.5
$ws.Range("X13").Value = 13
$ws.Range("AC13").Value = 7.5
$ws.Range("AF13").Value = 51
$ws.Range("AG13").Value = 401
$ws.Range("AL13").Value = 23
$ws.Range("AN13").Value = 4.75
$ws.Range("AT13").Value = 2.5
$ws.Range("AY13").Value = 15

